$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (pushes the old "tec da arquitetura" row
# and the trailing duplicate row down to 15 and 16 respectively). Copy the
# formatting from the row above first so the new row matches the table
# style (fill/border/number format + the 20pt row height).
$ws.Rows.Item(14).Insert()
$ws.Range("A13:J13").Copy()
$ws.Range("A14:J14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 20

# New row 14 data: "midia training" / Marketing course suggestion.
$ws.Cells.Item(14, 1).Value = "midia training "
$ws.Cells.Item(14, 2).Value = "Média"
$ws.Cells.Item(14, 3).Value = 12
$ws.Cells.Item(14, 4).Value = "Marketing"
$ws.Cells.Item(14, 5).Value = "Intermediário"
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = "Seminário"
$ws.Cells.Item(14, 8).Value = "brunocordeiro180"
$ws.Cells.Item(14, 9).Value = 88

# New "Aprovado" (approved) boolean column in J, matching the "Custo"
# column's style, with the header label and a flag per data row.
$ws.Range("I1:I16").Copy()
$ws.Range("J1:J16").PasteSpecial(-4122)

$ws.Cells.Item(1, 10).Value = "Aprovado"

$approvedByRow = @(
    @(2,  $false),
    @(3,  $true),
    @(4,  $false),
    @(5,  $true),
    @(6,  $true),
    @(7,  $false),
    @(8,  $false),
    @(9,  $false),
    @(10, $false),
    @(11, $false),
    @(12, $false),
    @(13, $false),
    @(14, $false),
    @(15, $false),
    @(16, $true)
)

foreach ($pair in $approvedByRow) {
    $ws.Cells.Item($pair[0], 10).Value = $pair[1]
}
